{"js": "// Update the worksheet date and the 26 division problems/answers that\n// changed between the \"before\" and \"after\" revisions of this document.\n// Each entry is [oldText, newText]; old texts are unique substrings in\n// the document (the date line plus 25 \"a\u00f7b=c, d\" table-cell strings),\n// so a simple search-and-replace-in-place is safe and unambiguous.\nconst replacements = [\n  [\"2025-04-22 Tuesday\", \"2025-04-23 Wednesday\"],\n  [\"45\u00f75=9, 0\", \"84\u00f73=28, 0\"],\n  [\"48\u00f77=6, 6\", \"23\u00f74=5, 3\"],\n  [\"86\u00f78=10, 6\", \"11\u00f74=2, 3\"],\n  [\"38\u00f79=4, 2\", \"52\u00f72=26, 0\"],\n  [\"62\u00f79=6, 8\", \"43\u00f76=7, 1\"],\n  [\"71\u00f72=35, 1\", \"22\u00f76=3, 4\"],\n  [\"83\u00f77=11, 6\", \"19\u00f72=9, 1\"],\n  [\"11\u00f72=5, 1\", \"78\u00f72=39, 0\"],\n  [\"35\u00f74=8, 3\", \"47\u00f73=15, 2\"],\n  [\"92\u00f79=10, 2\", \"63\u00f79=7, 0\"],\n  [\"10\u00f78=1, 2\", \"61\u00f79=6, 7\"],\n  [\"72\u00f72=36, 0\", \"90\u00f74=22, 2\"],\n  [\"32\u00f73=10, 2\", \"41\u00f73=13, 2\"],\n  [\"99\u00f77=14, 1\", \"23\u00f72=11, 1\"],\n  [\"67\u00f75=13, 2\", \"54\u00f76=9, 0\"],\n  [\"73\u00f73=24, 1\", \"95\u00f79=10, 5\"],\n  [\"68\u00f76=11, 2\", \"97\u00f79=10, 7\"],\n  [\"57\u00f72=28, 1\", \"76\u00f77=10, 6\"],\n  [\"38\u00f75=7, 3\", \"34\u00f77=4, 6\"],\n  [\"42\u00f74=10, 2\", \"74\u00f76=12, 2\"],\n  [\"19\u00f79=2, 1\", \"47\u00f78=5, 7\"],\n  [\"37\u00f75=7, 2\", \"14\u00f73=4, 2\"],\n  [\"59\u00f77=8, 3\", \"75\u00f78=9, 3\"],\n  [\"15\u00f78=1, 7\", \"99\u00f78=12, 3\"],\n  [\"28\u00f78=3, 4\", \"74\u00f74=18, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 26 division problems/answers that\n# changed between the \"before\" and \"after\" revisions of this document.\n# Each pair is (oldText, newText); old texts are unique substrings in the\n# document (the date line plus 25 \"a\u00f7b=c, d\" table-cell strings), so a\n# plain Find/Replace-all over the whole document body is unambiguous.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindWrap = 1  # wdFindContinue\n\n$replacements = @(\n    @(\"2025-04-22 Tuesday\", \"2025-04-23 Wednesday\"),\n    @(\"45\u00f75=9, 0\", \"84\u00f73=28, 0\"),\n    @(\"48\u00f77=6, 6\", \"23\u00f74=5, 3\"),\n    @(\"86\u00f78=10, 6\", \"11\u00f74=2, 3\"),\n    @(\"38\u00f79=4, 2\", \"52\u00f72=26, 0\"),\n    @(\"62\u00f79=6, 8\", \"43\u00f76=7, 1\"),\n    @(\"71\u00f72=35, 1\", \"22\u00f76=3, 4\"),\n    @(\"83\u00f77=11, 6\", \"19\u00f72=9, 1\"),\n    @(\"11\u00f72=5, 1\", \"78\u00f72=39, 0\"),\n    @(\"35\u00f74=8, 3\", \"47\u00f73=15, 2\"),\n    @(\"92\u00f79=10, 2\", \"63\u00f79=7, 0\"),\n    @(\"10\u00f78=1, 2\", \"61\u00f79=6, 7\"),\n    @(\"72\u00f72=36, 0\", \"90\u00f74=22, 2\"),\n    @(\"32\u00f73=10, 2\", \"41\u00f73=13, 2\"),\n    @(\"99\u00f77=14, 1\", \"23\u00f72=11, 1\"),\n    @(\"67\u00f75=13, 2\", \"54\u00f76=9, 0\"),\n    @(\"73\u00f73=24, 1\", \"95\u00f79=10, 5\"),\n    @(\"68\u00f76=11, 2\", \"97\u00f79=10, 7\"),\n    @(\"57\u00f72=28, 1\", \"76\u00f77=10, 6\"),\n    @(\"38\u00f75=7, 3\", \"34\u00f77=4, 6\"),\n    @(\"42\u00f74=10, 2\", \"74\u00f76=12, 2\"),\n    @(\"19\u00f79=2, 1\", \"47\u00f78=5, 7\"),\n    @(\"37\u00f75=7, 2\", \"14\u00f73=4, 2\"),\n    @(\"59\u00f77=8, 3\", \"75\u00f78=9, 3\"),\n    @(\"15\u00f78=1, 7\", \"99\u00f78=12, 3\"),\n    @(\"28\u00f78=3, 4\", \"74\u00f74=18, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindWrap, $false, $newText, $wdReplaceAll)\n}\n"}
